# edit.ps1 - applies the OOXML diff to before.pptx via PowerPoint COM automation
#
# Summary of the edit:
#   1. The "Fixed" date/time field shown on every slide (inherited from the
#      slide master + all 11 slide layouts) changes from 2021/10/8 to
#      2021/12/16.
#   2. On both content slides, two callout shapes get their red label text
#      extended:
#        "输入自变量列数" -> "输入指标变量" + "列数"
#        "输入因变量列数" -> "输入分组变量" + "列数"
#      (same visible final text, split across two runs - matching how the
#      text was extended/retyped in the source edit).

$p = $ppt.ActivePresentation

$oldDate = "2021/10/8"
$newDate = "2021/12/16"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# --- 1. Slide master date placeholder ---
$sm = $p.SlideMaster
Update-DatePlaceholder $sm.Shapes

# --- 2. Every slide layout's date placeholder ---
for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    Update-DatePlaceholder $cl.Shapes
}

# --- 3. Update the red callout labels on every slide ---
# Both callouts keep their original trailing "列数" characters and gain a
# new leading phrase; in the source edit this was typed/split slightly
# differently on the very last shape (new text inserted before the
# existing run instead of the existing run being retyped), so mirror that
# on the final shape (slide 2's "输入因变量列数" callout) while using the
# simpler "retype + append" approach everywhere else.
function Update-Callouts($slide, $useInsertBeforeForDep) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        if (-not $sh.TextFrame.HasText) { continue }
        $tr = $sh.TextFrame.TextRange
        $t = $tr.Text
        if ($t -eq "输入自变量列数") {
            $tr.Text = "输入指标变量"
            $tr.InsertAfter("列数")
        } elseif ($t -eq "输入因变量列数") {
            if ($useInsertBeforeForDep) {
                $tr.Characters(1, 5).Text = ""
                $tr.InsertBefore("输入分组变量")
            } else {
                $tr.Text = "输入分组变量"
                $tr.InsertAfter("列数")
            }
        }
    }
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $isLastSlide = ($si -eq $p.Slides.Count)
    Update-Callouts $p.Slides.Item($si) $isLastSlide
}
